$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New SVR parameter columns (K, L, M) ---
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# --- Clear the no-longer-present formatting on cells that used to carry
#     the duplicate "default" style (old cellXfs index 1) ---
$ws.Range("B1").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("H1").Style = "Normal"
$ws.Range("H2").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").Style = "Normal"
$ws.Range("A7").Style = "Normal"
$ws.Range("A8").Style = "Normal"
$ws.Range("A11").Style = "Normal"
$ws.Range("D14:G14").Style = "Normal"
$ws.Range("D15:G15").Style = "Normal"
$ws.Range("I15").Style = "Normal"
$ws.Range("D16:G16").Style = "Normal"
$ws.Range("I16").Style = "Normal"

# Row 13 only contained a styled-but-empty cell; clearing it removes the row entirely
$ws.Range("A13").Clear()

# --- Update the selection shown when the sheet is opened ---
$ws.Range("J6").Select()
